$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 02:34"

# --- Helper: locate a country's data row by its name in column A ---
$namesCol = $ws.Range("A4:A219")

function Set-CountryRow {
    param(
        [string]$CountryName,
        [double]$TotalCases,
        [double]$NewCases,
        [double]$ActiveCases,
        [double]$Recovered,
        [double]$DeathsToday,
        [double]$Deaths
    )
    $found = $namesCol.Find($CountryName)
    $r = $found.Row
    $ws.Cells.Item($r, 2).Value = $TotalCases
    $ws.Cells.Item($r, 3).Value = $NewCases
    $ws.Cells.Item($r, 4).Value = $ActiveCases
    $ws.Cells.Item($r, 5).Value = $Recovered
    $ws.Cells.Item($r, 7).Value = $DeathsToday
    $ws.Cells.Item($r, 8).Value = $Deaths
}

# --- Apply the updated per-country statistics ---
Set-CountryRow "Estados Unidos" 4861562 48010 2445854 2256801 546  158907
Set-CountryRow "Alemania"       212320  858   194700  8388   6    9232
Set-CountryRow "Argentina"      206743  4824  91302   111628 165  3813
Set-CountryRow "Panama"         68456   1003  42093   24866  26   1497
Set-CountryRow "Chequia"        17008   208   11708   4914   2    386
Set-CountryRow "Libia"          4063    226   625     3345   10   93
Set-CountryRow "Surinam"        1893    44    1227    639    0    27
Set-CountryRow "Bermudas"       157     0     144     4      0    9
Set-CountryRow "Seychelles"     114     0     113     1      0    0

# --- Re-sort the country table by "Casos totales" (column B) descending,
#     since Panama's update moves it above Kuwait / Bielorrusia ---
$sortRange = $ws.Range("A4:H219")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B4:B219"), 0, 2)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()
